# "update data with resort sheetname"
# The workbook has two sheets: "2022-Q2" (the per-fund holdings table) and
# "总计" (the quarter summary). This edit just re-sorts the sheet tabs so
# that the summary sheet "总计" comes first, followed by "2022-Q2".

$wb = $excel.ActiveWorkbook

$summarySheet = $wb.Worksheets.Item("总计")
$detailSheet  = $wb.Worksheets.Item("2022-Q2")

# Move "总计" so it sits right before "2022-Q2" -> becomes sheet 1,
# "2022-Q2" becomes sheet 2.
$summarySheet.Move($detailSheet)

# "总计" is now the first tab (position 1) and stays the active / selected
# sheet, matching the workbook's activeTab="0" selection state. Re-fetch a
# fresh reference by name so we activate the sheet at its new position.
$wb.Worksheets.Item("总计").Activate()
